# Natmi following Dr Hou advice
# Updates the LR-pairs (Vegfb-Flt1) sheet: the ligand- and receptor-expressing
# cell counts increase from 1 to 3, which changes the derived total/average
# expression and specificity statistics for every data row (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.911942
$ws.Range("H2").Value = 5.735825999999999
$ws.Range("I2").Value = 0.2156379149120961
$ws.Range("J2").Value = 0.2156379149120961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 149.829178
$ws.Range("N2").Value = 449.487534
$ws.Range("O2").Value = 0.965236887286734
$ws.Range("P2").Value = 0.965236887286734
$ws.Range("Q2").Value = 286.464698243676
$ws.Range("R2").Value = 2578.182284193084
$ws.Range("S2").Value = 0.2081416697707532
$ws.Range("T2").Value = 0.2081416697707532

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.911942
$ws.Range("H3").Value = 5.735825999999999
$ws.Range("I3").Value = 0.2156379149120961
$ws.Range("J3").Value = 0.2156379149120961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.021452666666667
$ws.Range("N3").Value = 3.064358
$ws.Range("O3").Value = 0.006580452523633729
$ws.Range("P3").Value = 0.006580452523633729
$ws.Range("Q3").Value = 1.952958254412
$ws.Range("R3").Value = 17.576624289708
$ws.Range("S3").Value = 0.001418995061374418
$ws.Range("T3").Value = 0.001418995061374418

# Row 4 (ECs -> sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.911942
$ws.Range("H4").Value = 5.735825999999999
$ws.Range("I4").Value = 0.2156379149120961
$ws.Range("J4").Value = 0.2156379149120961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.374661666666667
$ws.Range("N4").Value = 13.123985
$ws.Range("O4").Value = 0.02818266018963228
$ws.Range("P4").Value = 0.02818266018963228
$ws.Range("Q4").Value = 8.36409937629
$ws.Range("R4").Value = 75.27689438661
$ws.Range("S4").Value = 0.006077250079968444
$ws.Range("T4").Value = 0.006077250079968444

# Row 5 (FAPs -> ECs)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.898253666666667
$ws.Range("H5").Value = 8.694761
$ws.Range("I5").Value = 0.3268788371019294
$ws.Range("J5").Value = 0.3268788371019295
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 149.829178
$ws.Range("N5").Value = 449.487534
$ws.Range("O5").Value = 0.965236887286734
$ws.Range("P5").Value = 0.965236887286734
$ws.Range("Q5").Value = 434.2429645121527
$ws.Range("R5").Value = 3908.186680609374
$ws.Range("S5").Value = 0.3155155112441738
$ws.Range("T5").Value = 0.3155155112441738

# Row 6 (FAPs -> FAPs)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.898253666666667
$ws.Range("H6").Value = 8.694761
$ws.Range("I6").Value = 0.3268788371019294
$ws.Range("J6").Value = 0.3268788371019295
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.021452666666667
$ws.Range("N6").Value = 3.064358
$ws.Range("O6").Value = 0.006580452523633729
$ws.Range("P6").Value = 0.006580452523633729
$ws.Range("Q6").Value = 2.960428936493111
$ws.Range("R6").Value = 26.643860428438
$ws.Range("S6").Value = 0.00215101066852985
$ws.Range("T6").Value = 0.00215101066852985

# Row 7 (FAPs -> sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.898253666666667
$ws.Range("H7").Value = 8.694761
$ws.Range("I7").Value = 0.3268788371019294
$ws.Range("J7").Value = 0.3268788371019295
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.374661666666667
$ws.Range("N7").Value = 13.123985
$ws.Range("O7").Value = 0.02818266018963228
$ws.Range("P7").Value = 0.02818266018963228
$ws.Range("Q7").Value = 12.67887921584278
$ws.Range("R7").Value = 114.109912942585
$ws.Range("S7").Value = 0.009212315189225844
$ws.Range("T7").Value = 0.009212315189225844

# Row 8 (sCs -> ECs)
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.056250666666667
$ws.Range("H8").Value = 12.168752
$ws.Range("I8").Value = 0.4574832479859744
$ws.Range("J8").Value = 0.4574832479859745
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 149.829178
$ws.Range("N8").Value = 449.487534
$ws.Range("O8").Value = 0.965236887286734
$ws.Range("P8").Value = 0.965236887286734
$ws.Range("Q8").Value = 607.7447031486188
$ws.Range("R8").Value = 5469.702328337569
$ws.Range("S8").Value = 0.4415797062718069
$ws.Range("T8").Value = 0.441579706271807

# Row 9 (sCs -> FAPs)
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.056250666666667
$ws.Range("H9").Value = 12.168752
$ws.Range("I9").Value = 0.4574832479859744
$ws.Range("J9").Value = 0.4574832479859745
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.021452666666667
$ws.Range("N9").Value = 3.064358
$ws.Range("O9").Value = 0.006580452523633729
$ws.Range("P9").Value = 0.006580452523633729
$ws.Range("Q9").Value = 4.143268060135112
$ws.Range("R9").Value = 37.28941254121601
$ws.Range("S9").Value = 0.00301044679372946
$ws.Range("T9").Value = 0.003010446793729461

# Row 10 (sCs -> sCs)
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.056250666666667
$ws.Range("H10").Value = 12.168752
$ws.Range("I10").Value = 0.4574832479859744
$ws.Range("J10").Value = 0.4574832479859745
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.374661666666667
$ws.Range("N10").Value = 13.123985
$ws.Range("O10").Value = 0.02818266018963228
$ws.Range("P10").Value = 0.02818266018963228
$ws.Range("Q10").Value = 17.74472430185778
$ws.Range("R10").Value = 159.70251871672
$ws.Range("S10").Value = 0.01289309492043799
$ws.Range("T10").Value = 0.01289309492043799
